# Change the table style applied to the table on slide 5 (the only
# table in the deck) from "Table_0" ({2E6A23E8-7600-48ED-8ADC-3C4B186ECEDF})
# to the built-in style {713F5185-4D79-442A-977A-6EE6D58F67CE}.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{713F5185-4D79-442A-977A-6EE6D58F67CE}")
    }
}
